$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-20 23:18:20'
$ws.Range("O2").Value = '0.2 °C'
$ws.Range("E3").Value = '2026-02-20 23:18:22'
$ws.Range("E4").Value = '2026-02-20 23:18:25'
$ws.Range("H4").Value = '''61%'
$ws.Range("J4").Value = '1023.2 hPa'
$ws.Range("N4").Value = '3.1 °C 22:58 TU'
$ws.Range("O4").Value = '9.6 °C'
$ws.Range("E5").Value = '2026-02-20 23:18:28'
$ws.Range("H5").Value = '''91%'
$ws.Range("E6").Value = '2026-02-20 23:18:30'
$ws.Range("J6").Value = '1023.1 hPa'
$ws.Range("O6").Value = '9.2 °C'
$ws.Range("E7").Value = '2026-02-20 23:18:32'
$ws.Range("J7").Value = '1023.0 hPa'
$ws.Range("E8").Value = '2026-02-20 23:18:35'
$ws.Range("J8").Value = '1023.3 hPa'
$ws.Range("E9").Value = '2026-02-20 23:18:37'
$ws.Range("H9").Value = '''46%'
$ws.Range("O9").Value = '13.2 °C'
$ws.Range("E10").Value = '2026-02-20 23:18:40'
$ws.Range("H10").Value = '''81%'
$ws.Range("E11").Value = '2026-02-20 23:18:42'
$ws.Range("E12").Value = '2026-02-20 23:18:45'
$ws.Range("E13").Value = '2026-02-20 23:18:47'
$ws.Range("J13").Value = '1024.6 hPa'
$ws.Range("N13").Value = '-1.4 °C 22:35 TU'
$ws.Range("O13").Value = '5.7 °C'
$ws.Range("E14").Value = '2026-02-20 23:18:50'
$ws.Range("H14").Value = '''61%'
$ws.Range("O14").Value = '11.6 °C'
$ws.Range("E15").Value = '2026-02-20 23:18:52'
$ws.Range("E16").Value = '2026-02-20 23:18:54'
$ws.Range("E17").Value = '2026-02-20 23:18:56'
$ws.Range("M17").Value = '6.2 °C 22:55 TU'
$ws.Range("O17").Value = '3.4 °C'
$ws.Range("E18").Value = '2026-02-20 23:18:59'
$ws.Range("J18").Value = '1023.4 hPa'
$ws.Range("O18").Value = '7.5 °C'
$ws.Range("E19").Value = '2026-02-20 23:19:01'
$ws.Range("E20").Value = '2026-02-20 23:19:04'
$ws.Range("O20").Value = '-2.6 °C'
$ws.Range("E21").Value = '2026-02-20 23:19:06'
$ws.Range("H21").Value = '''39%'
$ws.Range("J21").Value = '1023.4 hPa'
$ws.Range("N21").Value = '2.4 °C 22:59 TU'
$ws.Range("O21").Value = '8.8 °C'
$ws.Range("E22").Value = '2026-02-20 23:19:09'
$ws.Range("E23").Value = '2026-02-20 23:19:11'
$ws.Range("O23").Value = '-4.4 °C'
$ws.Range("E24").Value = '2026-02-20 23:19:14'
$ws.Range("J24").Value = '1025.8 hPa'
$ws.Range("N24").Value = '4.4 °C 22:58 TU'
$ws.Range("O24").Value = '9.3 °C'
$ws.Range("E25").Value = '2026-02-20 23:19:16'
$ws.Range("E26").Value = '2026-02-20 23:19:18'
$ws.Range("H26").Value = '''36%'
$ws.Range("J26").Value = '1022.2 hPa'
$ws.Range("E27").Value = '2026-02-20 23:19:21'
$ws.Range("H27").Value = '''44%'
$ws.Range("O27").Value = '-0.6 °C'
$ws.Range("E28").Value = '2026-02-20 23:19:23'
$ws.Range("J28").Value = '1023.6 hPa'
$ws.Range("O28").Value = '6.7 °C'
$ws.Range("E29").Value = '2026-02-20 23:19:26'
$ws.Range("H29").Value = '''70%'
$ws.Range("O29").Value = '9.8 °C'
$ws.Range("E30").Value = '2026-02-20 23:19:28'
$ws.Range("J30").Value = '1022.9 hPa'
$ws.Range("O30").Value = '10.8 °C'
$ws.Range("E31").Value = '2026-02-20 23:19:30'
$ws.Range("J31").Value = '1022.1 hPa'
$ws.Range("E32").Value = '2026-02-20 23:19:33'
$ws.Range("K32").Value = '13.4 MJ/m2'
$ws.Range("E33").Value = '2026-02-20 23:19:35'
$ws.Range("H33").Value = '''44%'
$ws.Range("J33").Value = '1023.8 hPa'
$ws.Range("N33").Value = '1.5 °C 22:51 TU'
$ws.Range("O33").Value = '5.7 °C'
$ws.Range("E34").Value = '2026-02-20 23:19:38'
$ws.Range("E35").Value = '2026-02-20 23:19:40'
$ws.Range("H35").Value = '''77%'
$ws.Range("J35").Value = '1027.2 hPa'
$ws.Range("E36").Value = '2026-02-20 23:19:43'
$ws.Range("J36").Value = '1023.0 hPa'
$ws.Range("E37").Value = '2026-02-20 23:19:45'
$ws.Range("J37").Value = '1025.2 hPa'
$ws.Range("O37").Value = '4.3 °C'
$ws.Range("E38").Value = '2026-02-20 23:19:48'
$ws.Range("H38").Value = '''70%'
$ws.Range("O38").Value = '8.6 °C'
$ws.Range("E39").Value = '2026-02-20 23:19:50'
$ws.Range("O39").Value = '-2.3 °C'
$ws.Range("E40").Value = '2026-02-20 23:19:53'
$ws.Range("J40").Value = '1024.2 hPa'
$ws.Range("O40").Value = '9.9 °C'
$ws.Range("E41").Value = '2026-02-20 23:19:55'
$ws.Range("J41").Value = '1023.5 hPa'
$ws.Range("K41").Value = '14.1 MJ/m2'
$ws.Range("E42").Value = '2026-02-20 23:19:58'
$ws.Range("E43").Value = '2026-02-20 23:20:00'
$ws.Range("H43").Value = '''78%'
$ws.Range("O43").Value = '4.8 °C'
$ws.Range("E44").Value = '2026-02-20 23:20:02'
$ws.Range("M44").Value = '-0.1 °C 22:42 TU'
$ws.Range("O44").Value = '-4.2 °C'
$ws.Range("E45").Value = '2026-02-20 23:20:05'
$ws.Range("J45").Value = '1030.1 hPa'
$ws.Range("E46").Value = '2026-02-20 23:20:07'
$ws.Range("J46").Value = '1026.8 hPa'
$ws.Range("N46").Value = '7.8 °C 22:59 TU'
$ws.Range("O46").Value = '11.7 °C'
